# Apply the diff: rename "Gen" header to "MaxFES", rewrite the column A
# values (generation counts -> fraction-of-budget values), delete the
# "Run 50" column (AZ), and update the trailing "Mean" column values
# (now shifted from BA into AZ).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header rename: A1 "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 2. Column A data values (rows 2-14)
$colAValues = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $colAValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $colAValues[$i]
}

# 3. Delete the "Run 50" column (column AZ = 52). This shifts the old
#    "Mean" column (BA = 53) left into AZ, matching the diff's dimension
#    shrink from BA14 to AZ14.
$ws.Range("AZ:AZ").Delete()

# 4. Update the (now shifted) "Mean" column values in AZ2:AZ14
$meanValues = @(
    137.54410392,
    123.07174318,
    66.91472266,
    20.2295131,
    12.39226708,
    8.44348064,
    6.76509785,
    5.29377555,
    4.09015832,
    3.24828284,
    2.75511937,
    2.35238379,
    2.02443784
)
for ($i = 0; $i -lt $meanValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 52).Value = $meanValues[$i]
}
